$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "61+7="
$cell = $t.Cell(1, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "85-56="
$cell = $t.Cell(1, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "2+79="
$cell = $t.Cell(1, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "50-5="
$cell = $t.Cell(1, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "87-52="
$cell = $t.Cell(2, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "57+34="
$cell = $t.Cell(2, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "7+3="
$cell = $t.Cell(2, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "18+38="
$cell = $t.Cell(2, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "90-44="
$cell = $t.Cell(2, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "93-12="
$cell = $t.Cell(3, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "50+46="
$cell = $t.Cell(3, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "41+22="
$cell = $t.Cell(3, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "9+87="
$cell = $t.Cell(3, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "46-25="
$cell = $t.Cell(3, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "12+62="
$cell = $t.Cell(4, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "32-25="
$cell = $t.Cell(4, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "11+50="
$cell = $t.Cell(4, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "29+60="
$cell = $t.Cell(4, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "47-20="
$cell = $t.Cell(4, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "26+65="
$cell = $t.Cell(5, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "42-19="
$cell = $t.Cell(5, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "63+8="
$cell = $t.Cell(5, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "4+85="
$cell = $t.Cell(5, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "76-18="
$cell = $t.Cell(5, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "43-34="
$cell = $t.Cell(6, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "58+6="
$cell = $t.Cell(6, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "98-4="
$cell = $t.Cell(6, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "95-5="
$cell = $t.Cell(6, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "93-75="
$cell = $t.Cell(6, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "26-13="
$cell = $t.Cell(7, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "61+13="
$cell = $t.Cell(7, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "72+27="
$cell = $t.Cell(7, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "74-32="
$cell = $t.Cell(7, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "10+54="
$cell = $t.Cell(7, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "13+61="
$cell = $t.Cell(8, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "60+37="
$cell = $t.Cell(8, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "57+30="
$cell = $t.Cell(8, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "91-14="
$cell = $t.Cell(8, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "59-47="
$cell = $t.Cell(8, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "46+25="
$cell = $t.Cell(9, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "27+71="
$cell = $t.Cell(9, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "74-42="
$cell = $t.Cell(9, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "59-17="
$cell = $t.Cell(9, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "3+48="
$cell = $t.Cell(9, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "9+67="
$cell = $t.Cell(10, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "93-20="
$cell = $t.Cell(10, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "29-4="
$cell = $t.Cell(10, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "37+41="
$cell = $t.Cell(10, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "90-30="
$cell = $t.Cell(10, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "32-3="
$cell = $t.Cell(11, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "4+20="
$cell = $t.Cell(11, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "78-77="
$cell = $t.Cell(11, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "75+23="
$cell = $t.Cell(11, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "26-12="
$cell = $t.Cell(11, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "17+5="
$cell = $t.Cell(12, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "50-43="
$cell = $t.Cell(12, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "87+8="
$cell = $t.Cell(12, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "56+14="
$cell = $t.Cell(12, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "1+85="
$cell = $t.Cell(12, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "8-2="
$cell = $t.Cell(13, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "75+2="
$cell = $t.Cell(13, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "48+18="
$cell = $t.Cell(13, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "37+16="
$cell = $t.Cell(13, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "28+19="
$cell = $t.Cell(13, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "80-33="
$cell = $t.Cell(14, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "86-29="
$cell = $t.Cell(14, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "82-55="
$cell = $t.Cell(14, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "68+17="
$cell = $t.Cell(14, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "78-35="
$cell = $t.Cell(14, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "91-49="
$cell = $t.Cell(15, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "25+53="
$cell = $t.Cell(15, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "50-47="
$cell = $t.Cell(15, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "59-20="
$cell = $t.Cell(15, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "42-33="
$cell = $t.Cell(15, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "79-67="
$cell = $t.Cell(16, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "79-61="
$cell = $t.Cell(16, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "57+39="
$cell = $t.Cell(16, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "64-22="
$cell = $t.Cell(16, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "66-35="
$cell = $t.Cell(16, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "18-9="
$cell = $t.Cell(17, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "89-63="
$cell = $t.Cell(17, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "70-22="
$cell = $t.Cell(17, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "30+47="
$cell = $t.Cell(17, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "4+87="
$cell = $t.Cell(17, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "95-48="
$cell = $t.Cell(18, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "14+29="
$cell = $t.Cell(18, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "70-63="
$cell = $t.Cell(18, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "2+87="
$cell = $t.Cell(18, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "65-12="
$cell = $t.Cell(18, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "4+70="
$cell = $t.Cell(19, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "51-35="
$cell = $t.Cell(19, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "62+11="
$cell = $t.Cell(19, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "51+30="
$cell = $t.Cell(19, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "80-46="
$cell = $t.Cell(19, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "3+42="
$cell = $t.Cell(20, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "70-50="
$cell = $t.Cell(20, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "62+6="
$cell = $t.Cell(20, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "54+3="
$cell = $t.Cell(20, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "90-52="
$cell = $t.Cell(20, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "77-46="
